$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.41

$ws.Range("B3").Value = 1.41
$ws.Range("D3").Value = 1.35
$ws.Range("E3").Value = 1.37
$ws.Range("F3").Value = 1.11

$ws.Range("C4").Value = 1.5
$ws.Range("D4").Value = 1.38
$ws.Range("E4").Value = 1.2

$ws.Range("C5").Value = 1.29
$ws.Range("D5").Value = 1.35
$ws.Range("F5").Value = 1.07
$ws.Range("G5").Value = 0.52

$ws.Range("D6").Value = 1.61
$ws.Range("E6").Value = 1.31

$ws.Range("E7").Value = 2.14
$ws.Range("F7").Value = 1.45
